$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. The sheet that currently carries tabSelected ("TUC - FWYNTK Submenu
#    Names") stops being the active tab; it keeps a plain cell selection.
# ---------------------------------------------------------------------------
$wwynk = $wb.Worksheets.Item("TUC - FWYNTK Submenu Names")
$wwynk.Range("Q15").Select()

# ---------------------------------------------------------------------------
# 2. Add the 7 new "CUC" sheets, in final tab order, at the end of the
#    workbook.
# ---------------------------------------------------------------------------
$sheetNames = @(
  "CUC - Grid Header Names",
  "CUC - TravelFlexibility Names",
  "CUC - TravelFlexibility URLs",
  "CUC - TravelingWithUs Names",
  "CUC - TravelingWithUs URLs",
  "CUC - CaringForYou Names",
  "CUC - CaringForYou URLs"
)
foreach ($name in $sheetNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $s = $wb.Worksheets.Add($null, $lastSheet)
    $s.Name = $name
}

$gridHeaders    = $wb.Worksheets.Item("CUC - Grid Header Names")
$flexNames      = $wb.Worksheets.Item("CUC - TravelFlexibility Names")
$flexUrls       = $wb.Worksheets.Item("CUC - TravelFlexibility URLs")
$travelingNames = $wb.Worksheets.Item("CUC - TravelingWithUs Names")
$travelingUrls  = $wb.Worksheets.Item("CUC - TravelingWithUs URLs")
$caringNames    = $wb.Worksheets.Item("CUC - CaringForYou Names")
$caringUrls     = $wb.Worksheets.Item("CUC - CaringForYou URLs")

# ---------------------------------------------------------------------------
# 3. Populate values in the same order the original author typed them in,
#    so the shared-string table comes out in the same sequence.
# ---------------------------------------------------------------------------

# Grid headers
$gridHeaders.Range("A1").Value = "TRAVEL FLEXIBILITY"
$gridHeaders.Range("A2").Value = "TRAVELING WITH US"
$gridHeaders.Range("A3").Value = "CARING FOR YOU"

# TravelFlexibility URLs
$flexUrls.Range("A1").Value = "https://www.delta.com/us/en/travel-update-center/overview#waiver"
$flexUrls.Range("A2").Value = "https://www.delta.com/us/en/travel-update-center/overview#confidence"
$flexUrls.Range("A3").Value = "https://www.delta.com/us/en/travel-update-center/overview#faq"

# TravelingWithUs URLs
$travelingUrls.Range("A1").Value = "https://www.delta.com/us/en/travel-update-center/overview#deltaclean"
$travelingUrls.Range("A2").Value = "https://www.delta.com/us/en/travel-update-center/overview#skyclub"
$travelingUrls.Range("A3").Value = "https://www.delta.com/us/en/travel-update-center/overview#flydeltaapp"

# CaringForYou URLs (note: A3 repeats the same URL as A2, by design)
$caringUrls.Range("A1").Value = "https://www.delta.com/us/en/travel-update-center/overview#skymiles"
$caringUrls.Range("A2").Value = "https://www.delta.com/us/en/travel-update-center/overview#frontlines"
$caringUrls.Range("A3").Value = "https://www.delta.com/us/en/travel-update-center/overview#frontlines"

# TravelFlexibility Names
$flexNames.Range("A1").Value = "Updates, Waivers and eCredits`n, Go to footer note"
$flexNames.Range("A2").Value = "Booking with Confidence`n, Go to footer note"
$flexNames.Range("A3").Value = "Frequently Asked Questions`n, Go to footer note"

# TravelingWithUs Names - A1 first
$travelingNames.Range("A1").Value = "Standard for Safer Travel`n, Go to footer note"

# CaringForYou Names
$caringNames.Range("A1").Value = "SkyMiles® Program Updates`n, Go to footer note"
$caringNames.Range("A2").Value = "Supporting Medical Volunteers`n, Go to footer note"
$caringNames.Range("A3").Value = "Protective Equipment for Healthcare Workers`n, Go to footer note"

# TravelingWithUs Names - A2/A3 afterwards
$travelingNames.Range("A2").Value = "Delta Sky Club Updates`n, Go to footer note"
$travelingNames.Range("A3").Value = "Download the Fly Delta App`n, Go to footer note"

# ---------------------------------------------------------------------------
# 4. Formatting: wrap text + 30pt row height on the three "Names" sheets,
#    column widths (best-fit approximations) everywhere, and the
#    per-sheet cell selections recorded in the workbook.
# ---------------------------------------------------------------------------

foreach ($ws in @($flexNames, $travelingNames, $caringNames)) {
    $ws.Range("A1:A3").WrapText = $true
    $ws.Rows.Item(1).RowHeight = 30
    $ws.Rows.Item(2).RowHeight = 30
    $ws.Rows.Item(3).RowHeight = 30
}

$flexNames.Columns.Item(1).ColumnWidth      = 28.5703125
$flexUrls.Columns.Item(1).ColumnWidth       = 69.42578125
$travelingNames.Columns.Item(1).ColumnWidth = 26.140625
$travelingUrls.Columns.Item(1).ColumnWidth  = 68.85546875
$caringNames.Columns.Item(1).ColumnWidth    = 42.28515625
$caringUrls.Columns.Item(1).ColumnWidth     = 67.28515625

$flexNames.PageSetup.Orientation = 1
$flexNames.Range("D7:D8").Select()
$flexUrls.Range("A3").Select()
$travelingUrls.Range("A3").Select()
$caringNames.Range("D6:D7").Select()
$caringUrls.Range("A3").Select()

# ---------------------------------------------------------------------------
# 5. "CUC - TravelingWithUs Names" ends up the active sheet/tab, with A4
#    selected; scroll the tab strip so the later tabs are in view.
# ---------------------------------------------------------------------------
$travelingNames.Activate()
$travelingNames.Range("A4").Select()
$excel.ActiveWindow.ScrollWorkbookTabs(8)
